$wb = $excel.ActiveWorkbook

# --- systemComponent sheet -------------------------------------------------
$ws2 = $wb.Worksheets.Item("systemComponent")

# Update the "notes" column header wording
$ws2.Range("F1").Value = "Brief description of the functionality that the component contributes to the system"

# Narrow column E (SyRO role) width
$ws2.Columns.Item(5).ColumnWidth = 11.998697916666666

# Move the cursor/selection to F3
$ws2.Activate() | Out-Null
$ws2.Range("F3").Select() | Out-Null

# --- component sheet --------------------------------------------------------
$ws3 = $wb.Worksheets.Item("component")

# Add a new "notes" column (G) to the component table
$ws3.Columns.Item(7).ColumnWidth = 25.666666666666668
$ws3.Range("G1").Value = "Brief description of  properties that are  important to the component regardless of its role in the system"
$ws3.Range("G2").Value = "notes"

# Header row grew taller to fit the new wrapped header text
$ws3.Rows.Item(1).RowHeight = 36

# component becomes the active sheet/tab, with G1 selected
$ws3.Activate() | Out-Null
$ws3.Range("G1").Select() | Out-Null
